# Lattice multiplication exercises: refresh all 15 practice-cell values
# (same 5x3 table shape; only the multiplicands/partial-product digits change).
$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$newCellText = @(
    @(@('62 x 18','  1    8','  ----','6|    |','2|    |'), @('52 x 15','  1    5','  ----','5|    |','2|    |'), @('91 x 65','  6    5','  ----','9|    |','1|    |')),
    @(@('34 x 12','  1    2','  ----','3|    |','4|    |'), @('18 x 66','  6    6','  ----','1|    |','8|    |'), @('55 x 11','  1    1','  ----','5|    |','5|    |')),
    @(@('85 x 30','  3    0','  ----','8|    |','5|    |'), @('58 x 35','  3    5','  ----','5|    |','8|    |'), @('84 x 22','  2    2','  ----','8|    |','4|    |')),
    @(@('97 x 67','  6    7','  ----','9|    |','7|    |'), @('44 x 50','  5    0','  ----','4|    |','4|    |'), @('84 x 68','  6    8','  ----','8|    |','4|    |')),
    @(@('40 x 90','  9    0','  ----','4|    |','0|    |'), @('48 x 54','  5    4','  ----','4|    |','8|    |'), @('24 x 79','  7    9','  ----','2|    |','4|    |'))
)

for ($ri = 1; $ri -le $tbl.Rows.Count; $ri++) {
  for ($ci = 1; $ci -le $tbl.Rows.Item($ri).Cells.Count; $ci++) {
    $lines = $newCellText[$ri - 1][$ci - 1]
    $runXml = '<w:r><w:rPr><w:sz w:val="32"/></w:rPr>'
    for ($li = 0; $li -lt $lines.Length; $li++) {
      if ($li -gt 0) { $runXml += "<w:br/>" }
      $line = $lines[$li]
      $escaped = $line.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
      if ($line -ne $line.Trim() -or $line -eq "") {
        $runXml += ('<w:t xml:space="preserve">' + $escaped + "</w:t>")
      } else {
        $runXml += ("<w:t>" + $escaped + "</w:t>")
      }
    }
    $runXml += "</w:r>"
    $pXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $runXml + "</w:p>"

    $cell = $tbl.Rows.Item($ri).Cells.Item($ci)
    $cellRange = $cell.Range
    # Cell.Range.End sits one position past the final paragraph mark (into the
    # cell-end mark); trimming that last position makes InsertXML replace the
    # whole paragraph instead of just appending a sibling one.
    $target = $d.Range($cellRange.Start, $cellRange.End - 1)
    $target.InsertXML($pXml) | Out-Null
  }
}

Write-Output "Lattice multiplication exercise values updated."